# manual dislocation uploading 2021/07/26 20:00
#
# Re-keys the "hours" (column B) values on Sheet1 of the KBL workbook for a
# number of dates, and restores the window/sheet view state (scroll
# position + selection) to where the user had navigated when saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes: column B values were re-keyed for several dates -------
$bChanges = @{
    "B2"  = 48
    "B5"  = 48
    "B7"  = 0
    "B8"  = 48
    "B9"  = 0
    "B10" = 48
    "B12" = 48
    "B14" = 48
    "B16" = 48
    "B18" = 48
    "B20" = 48
    "B23" = 0
    "B25" = 0
    "B27" = 0
    "B29" = 0
    "B31" = 0
}

foreach ($addr in $bChanges.Keys) {
    $ws.Range($addr).Value = $bChanges[$addr]
}

# --- View state: the sheet had been scrolled down (top row ~15) and the
#     active selection moved to B31 by the time the workbook was saved.
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1

# --- The workbook window had also moved down slightly on screen.
$excel.ActiveWindow.Top = 7800
